# Rename the inline-picture "name" metadata (wp:docPr / InlineShape.Name)
# for the Pearson logo (in both footers) and the BTec logo (in both
# headers), matching the picture back up with the correct imageN.ext
# label following the media re-numbering:
#   Pearson logo:  image2.png -> image1.png   (footers, default + first page)
#   BTec logo:     image1.jpg -> image2.jpg   (headers, default + first page)
#
# NOTE: direct $range.InlineShapes(1).Name = "..." assignments are not
# reliably committed when the InlineShape lives in a Footer story in this
# host, so the shape is selected first and renamed through
# $word.Selection.InlineShapes(1), which persists for both Headers and
# Footers.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Rename-InlineLogo($range, $newName) {
    $ishp = $range.InlineShapes(1)
    $ishp.Select()
    $selShape = $word.Selection.InlineShapes(1)
    $selShape.Name = $newName
}

# Footers (Pearson Edexcel logo): image2.png -> image1.png
Rename-InlineLogo $sec.Footers(1).Range "image1.png"
Rename-InlineLogo $sec.Footers(2).Range "image1.png"

# Headers (BTEC logo): image1.jpg -> image2.jpg
Rename-InlineLogo $sec.Headers(1).Range "image2.jpg"
Rename-InlineLogo $sec.Headers(2).Range "image2.jpg"

Write-Output "Renamed header/footer logo InlineShapes"
